$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift values for C1, D1, E1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows 2-5: column C becomes the genus text (same as D),
# column D stays the same genus text, column E becomes numeric 1
for ($r = 2; $r -le 5; $r++) {
    $genus = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $genus
    $ws.Cells.Item($r, 5).Value = 1
}
